$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 321
$range = $ws.Range("C2:C$lastRow")
$range.Value = 46082
